# The workbook was re-opened/re-saved (from a Windows machine's Excel
# into a Mac-hosted Excel build), which is why most of this diff is just
# incidental re-save metadata (fileVersion/rupBuild, the absPath used for
# the x15ac:absPath hint, window geometry, and the x14ac:dyDescent /
# autofit-derived row-height & column-width jitter). None of that is
# real "data" and none of it is exposed through the Excel object model in
# a way a user script could drive deliberately — it is whatever the host
# application stamps on save.
#
# The one deliberate, user-visible change captured in the diff is the
# worksheet being renamed from "T_dis" to "sheet1". Everything else
# (every <v> cell value in the sheet) is byte-for-byte identical before
# and after, which confirms this commit did not touch any data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "sheet1"
